# Apply updated daily price records (weekly consolidation) to rows 2-39.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; D=44350; I="Primera"; J=100; K=13000; L=14000; M=13500; P=750 },
    @{ Row=3; D=44350; I="Segunda"; J=50; K=11000; L=11000; M=11000; P=611 },
    @{ Row=4; D=44356; I="Primera"; J=100; K=10000; L=11000; M=10500; P=583 },
    @{ Row=5; D=44356; I="Segunda"; J=50; K=9000; L=9000; M=9000; P=500 },
    @{ Row=6; D=44349; I="Primera"; J=200; K=11000; L=12000; M=11500; P=639 },
    @{ Row=7; D=44349; I="Segunda"; J=100; K=10000; L=10000; M=10000; P=556 },
    @{ Row=8; D=44280; I="Primera"; J=200; K=10000; L=11000; M=10500; P=583 },
    @{ Row=9; D=44280; I="Segunda"; J=100; K=9000; L=9000; M=9000; P=500 },
    @{ Row=10; D=44364; I="Primera"; J=100; K=13000; L=14000; M=13500; P=750 },
    @{ Row=11; D=44364; I="Segunda"; J=50; K=11000; L=11000; M=11000; P=611 },
    @{ Row=12; D=44308; I="Primera"; J=200; K=10000; L=11000; M=10500; P=583 },
    @{ Row=13; D=44308; I="Segunda"; J=50; K=8000; L=8000; M=8000; P=444 },
    @{ Row=14; D=44328; I="Primera"; J=200; K=9000; L=10000; M=9500; P=528 },
    @{ Row=15; D=44328; I="Segunda"; J=100; K=8000; L=8000; M=8000; P=444 },
    @{ Row=16; D=44320; I="Primera"; J=100; K=9000; L=10000; M=9500; P=528 },
    @{ Row=17; D=44320; I="Segunda"; J=50; K=8000; L=8000; M=8000; P=444 },
    @{ Row=18; D=44265; I="Primera"; J=100; K=13000; L=14000; M=13500; P=750 },
    @{ Row=19; D=44259; I="Primera"; J=100; K=12000; L=13000; M=12500; P=694 },
    @{ Row=20; D=44259; I="Segunda"; J=50; K=10000; L=10000; M=10000; P=556 },
    @{ Row=21; D=44392; I="Primera"; J=200; K=15000; L=16000; M=15500; P=861 },
    @{ Row=22; D=44392; I="Segunda"; J=100; K=14000; L=14000; M=14000; P=778 },
    @{ Row=23; D=44384; I="Primera"; J=200; K=15000; L=16000; M=15500; P=861 },
    @{ Row=24; D=44384; I="Segunda"; J=50; K=13000; L=13000; M=13000; P=722 },
    @{ Row=25; D=44272; I="Primera"; J=100; K=10000; L=11000; M=10500; P=583 },
    @{ Row=26; D=44272; I="Segunda"; J=50; K=9000; L=9000; M=9000; P=500 },
    @{ Row=27; D=44397; I="Primera"; J=100; K=14000; L=15000; M=14500; P=806 },
    @{ Row=28; D=44335; I="Primera"; J=100; K=12000; L=13000; M=12500; P=694 },
    @{ Row=29; D=44335; I="Segunda"; J=50; K=10000; L=10000; M=10000; P=556 },
    @{ Row=30; D=44316; I="Primera"; J=200; K=10000; L=11000; M=10500; P=583 },
    @{ Row=31; D=44316; I="Segunda"; J=100; K=9000; L=9000; M=9000; P=500 },
    @{ Row=32; D=44253; I="Primera"; J=100; K=12000; L=12000; M=12000; P=667 },
    @{ Row=33; D=44253; I="Segunda"; J=100; K=10000; L=10000; M=10000; P=556 },
    @{ Row=34; D=44342; I="Primera"; J=100; K=11000; L=12000; M=11500; P=639 },
    @{ Row=35; D=44342; I="Segunda"; J=50; K=9000; L=9000; M=9000; P=500 },
    @{ Row=36; D=44379; I="Primera"; J=200; K=15000; L=16000; M=15500; P=861 },
    @{ Row=37; D=44379; I="Segunda"; J=100; K=13000; L=13000; M=13000; P=722 },
    @{ Row=38; D=44313; I="Primera"; J=200; K=10000; L=11000; M=10500; P=583 },
    @{ Row=39; D=44313; I="Segunda"; J=100; K=9000; L=9000; M=9000; P=500 }
)

foreach ($r in $rows) {
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("I" + $r.Row).Value = $r.I
    $ws.Range("J" + $r.Row).Value = $r.J
    $ws.Range("K" + $r.Row).Value = $r.K
    $ws.Range("L" + $r.Row).Value = $r.L
    $ws.Range("M" + $r.Row).Value = $r.M
    $ws.Range("P" + $r.Row).Value = $r.P
}